# RAD Phase 3 Summary and BeforePayments test data update:
# refresh the recorded run timestamps in column B (Date) for rows 2 and 4.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Sun Jan 14 16:34:27 EST 2024"
$ws.Range("B4").Value = "Sun Jan 14 16:34:39 EST 2024"
